# Insert a new weekly record at row 59 (pushing existing rows 59-127 down to 60-128)
# and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 59:127 down by inserting a new row at 59.
$ws.Rows.Item(59).Insert()

# Fixed/common columns (identical across all data rows in this sheet).
$ws.Cells.Item(59, 1).Value = 8
$ws.Cells.Item(59, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(59, 3).Value = "Coquimbo"
$ws.Cells.Item(59, 4).Value = 45049
$ws.Cells.Item(59, 5).Value = 4
$ws.Cells.Item(59, 6).Value = 100114007
$ws.Cells.Item(59, 7).Value = "Jengibre"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 500
$ws.Cells.Item(59, 11).Value = 17000
$ws.Cells.Item(59, 12).Value = 18000
$ws.Cells.Item(59, 13).Value = 17500
$ws.Cells.Item(59, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(59, 15).Value = "Perú"
$ws.Cells.Item(59, 16).Value = 1346
$ws.Cells.Item(59, 17).Value = 13
$ws.Cells.Item(59, 18).Value = "Hortaliza"

# Match the date-number-format style used by the rest of column D.
$ws.Cells.Item(59, 4).NumberFormat = $ws.Cells.Item(60, 4).NumberFormat
